# Ofertas.xlsx — "list_offers" sheet update
#
# 1. Refresh the tracking `srsltid` query-string parameter on two existing
#    offer links (rows 9 and 11) — the product rows themselves are unchanged.
# 2. Insert a brand-new offer ("ailos aproxima | iphone 12 64gb azul - swap")
#    above the existing "smartphone apple iphone 12 64gb câmera dupla" row,
#    which pushes it (with an updated price) and the row below it
#    ("placa de video ... rtx 3060 ti ...") down by one row.
# 3. Re-stripe the shifted rows so the banded (style 2 / style 3) pattern
#    keeps alternating correctly through the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the srsltid query parameter on the two existing URLs -------
$ws.Range("C9").Value = "https://www.horizonplay.com.br/apple/iphone/apple-iphone-12-64gb-azul-novo-lacrado-tela-super-retina-xdr-oled-6-1?variant_id=21019&parceiro=8926&srsltid=AfmBOopqbiO0i-O8ok7RtRlBKV1bYVWIvzNe224FJHUrq61-q-tWbwQ2CIU"
$ws.Range("C11").Value = "https://lumixpel.lojavirtualnuvem.com.br/produtos/iphone-11-apple-64gb-e-128gb-preto-61-12mp-ios/?variant=660570516&pf=mc&srsltid=AfmBOop3hX52tE7OjosJ4QLXIMq4G1OmpC0TR4CEadcctajK-hJ_gyFcm70"

# --- 2) Insert the new offer row above row 12 ------------------------------
# This shifts the old row 12 ("smartphone ...") to row 13 and the old row 13
# ("placa de video ...") to row 14, carrying their formatting down with them.
$ws.Rows.Item(12).Insert()

# --- 3) Re-stripe rows 12-14 so the alternating banding keeps matching ----
# Row 12 (new) takes the "even" banding used by rows 2,4,6,8,10.
$ws.Range("A10:C10").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 13 (shifted-down old row 12) takes the "odd" banding used by rows 3,5,7,9,11.
$ws.Range("A11:C11").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 14 (shifted-down old row 13) takes the "even" banding.
$ws.Range("A10:C10").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4) Fill in the new offer's data (row 12) ------------------------------
$ws.Range("A12").Value = "ailos aproxima | iphone 12 64gb azul - swap"
$ws.Range("B12").Value = 3424
$ws.Range("C12").Value = "https://ailosaproxima.coop.br/loja/malibu-shop/produto/380030/iphone-12-64gb-azul-swap"

# --- 5) Update the price on the shifted "smartphone" offer (now row 13) ---
$ws.Range("B13").Value = 3039
